# Correct csv module error handling
# - Append 3 new rows (62-64) of parsed CSV data to each of the 4 sheets.
# - On sheets DE_PLT_#1 and DE_PLT_#2, correct previously mis-decoded
#   "actual length" values (column D / D_DEC column H) for a handful of
#   existing rows so that the decoded byte value matches the checksum-derived
#   decimal value.

function Set-RowData {
    param($ws, $row, $a, $b, $c, $d, $e, $f, $g, $h, $i)
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)
$sheet3 = $wb.Worksheets.Item(3)
$sheet4 = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------------
# Sheet 1 (DE_LFT_#1): fix existing rows -- none; just append rows 62-64
# ---------------------------------------------------------------------------
$g1 = [double]"7.598631275147109e+23"
Set-RowData $sheet1 62 45848.43368055556 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x58" "0x14" 380 $g1 344 14
Set-RowData $sheet1 63 45849.43524305556 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x58" "0x14" 380 $g1 344 14
Set-RowData $sheet1 64 45850.43918981482 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x54" "0x14" 380 $g1 340 14

# ---------------------------------------------------------------------------
# Sheet 2 (DE_LFT_#2): append rows 62-64
# ---------------------------------------------------------------------------
$g2 = [double]"5.68432987514711e+23"
Set-RowData $sheet2 62 45848.43368055556 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x58" "0xe" 380 $g2 344 14
Set-RowData $sheet2 63 45849.43524305556 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x58" "0xe" 380 $g2 344 14
Set-RowData $sheet2 64 45850.43918981482 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x54" "0xe" 380 $g2 340 14

# ---------------------------------------------------------------------------
# Sheet 3 (DE_PLT_#1): correct D/H for rows 28-34, 47-50, 56-60, then append
# ---------------------------------------------------------------------------
$sheet3.Range("D28").Value = "0x00,0x81"
$sheet3.Range("D29").Value = "0x00,0x81"
$sheet3.Range("D30").Value = "0x00,0x81"
$sheet3.Range("D31").Value = "0x00,0x81"
$sheet3.Range("D32").Value = "0x00,0x81"
$sheet3.Range("D33").Value = "0x00,0x81"

$sheet3.Range("D34").Value = "0x00,0x81"
$sheet3.Range("H34").Value = 129

$sheet3.Range("D47").Value = "0x00,0x7E"
$sheet3.Range("H47").Value = 126
$sheet3.Range("D48").Value = "0x00,0x7E"
$sheet3.Range("H48").Value = 126
$sheet3.Range("D49").Value = "0x00,0x7E"
$sheet3.Range("H49").Value = 126
$sheet3.Range("D50").Value = "0x00,0x7E"
$sheet3.Range("H50").Value = 126

$sheet3.Range("D56").Value = "0x00,0x7D"
$sheet3.Range("H56").Value = 125
$sheet3.Range("D57").Value = "0x00,0x7D"
$sheet3.Range("H57").Value = 125
$sheet3.Range("D58").Value = "0x00,0x7D"
$sheet3.Range("H58").Value = 125
$sheet3.Range("D59").Value = "0x00,0x7D"
$sheet3.Range("H59").Value = 125
$sheet3.Range("D60").Value = "0x00,0x7D"
$sheet3.Range("H60").Value = 125

$g3 = [double]"5.68631262647114e+23"
Set-RowData $sheet3 62 45848.43368055556 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x7C" "0x7" 130 $g3 124 7
Set-RowData $sheet3 63 45849.43524305556 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x7C" "0x7" 130 $g3 124 7
Set-RowData $sheet3 64 45850.43918981482 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x7C" "0x7" 130 $g3 124 7

# ---------------------------------------------------------------------------
# Sheet 4 (DE_PLT_#2): correct D/H for rows 28-34, 47-50, 56-60, then append
# ---------------------------------------------------------------------------
$sheet4.Range("D28").Value = "0x00,0x82"
$sheet4.Range("D29").Value = "0x00,0x82"
$sheet4.Range("D30").Value = "0x00,0x82"
$sheet4.Range("D31").Value = "0x00,0x82"
$sheet4.Range("D32").Value = "0x00,0x82"
$sheet4.Range("D33").Value = "0x00,0x82"

$sheet4.Range("D34").Value = "0x00,0x80"
$sheet4.Range("H34").Value = 128

$sheet4.Range("D47").Value = "0x00,0x7F"
$sheet4.Range("H47").Value = 127
$sheet4.Range("D48").Value = "0x00,0x7F"
$sheet4.Range("H48").Value = 127
$sheet4.Range("D49").Value = "0x00,0x7F"
$sheet4.Range("H49").Value = 127
$sheet4.Range("D50").Value = "0x00,0x7F"
$sheet4.Range("H50").Value = 127

$sheet4.Range("D56").Value = "0x00,0x7E"
$sheet4.Range("H56").Value = 126
$sheet4.Range("D57").Value = "0x00,0x7E"
$sheet4.Range("H57").Value = 126
$sheet4.Range("D58").Value = "0x00,0x7E"
$sheet4.Range("H58").Value = 126
$sheet4.Range("D59").Value = "0x00,0x7E"
$sheet4.Range("H59").Value = 126
$sheet4.Range("D60").Value = "0x00,0x7E"
$sheet4.Range("H60").Value = 126

$g4 = [double]"9.85046333984776e+23"
Set-RowData $sheet4 62 45848.43368055556 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x7C" "0x3" 130 $g4 124 3
Set-RowData $sheet4 63 45849.43524305556 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x7C" "0x3" 130 $g4 124 3
Set-RowData $sheet4 64 45850.43918981482 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x7C" "0x3" 130 $g4 124 3
